# Apply cell updates described by the commit diff (cryptos list refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.108.41"
$ws.Range("E2").Value = "  -3.04%  "
$ws.Range("D3").Value = "1.908.47"
$ws.Range("E3").Value = "  -3.81%  "
$ws.Range("D4").Value = "'1.003"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -1.24%  "
$ws.Range("D5").Value = "'327.73"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.48%  "
$ws.Range("E6").Value = "  -1.29%  "
$ws.Range("D7").Value = "'0.4642"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -5.38%  "
$ws.Range("D8").Value = "'0.4004"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -3.49%  "
$ws.Range("D9").Value = "'53.16"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -3.36%  "
$ws.Range("D10").Value = "'0.08378"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -5.01%  "
$ws.Range("D11").Value = "'1.041"
$ws.Range("D11").ClearFormats()
$ws.Range("D12").Value = "'21.92"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -3.34%  "
$ws.Range("D13").Value = "1.905.70"
$ws.Range("E13").Value = "  -8.44%  "
$ws.Range("D14").Value = "'7.402"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -5.84%  "
$ws.Range("D15").Value = "'6.043"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -4.81%  "
$ws.Range("D16").Value = "'1.005"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -1.18%  "
$ws.Range("D17").Value = "'89.37"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -2.11%  "
$ws.Range("D18").Value = "'0.00001065"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -2.78%  "
$ws.Range("D19").Value = "'0.06602"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -0.97%  "
$ws.Range("D20").Value = "'17.84"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -7.40%  "
$ws.Range("D21").Value = "'1.002"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -1.30%  "
$ws.Range("D22").Value = "'5.730"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -3.35%  "
$ws.Range("D23").Value = "28.113.96"
$ws.Range("E23").Value = "  -3.18%  "
$ws.Range("D24").Value = "'11.16"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -5.57%  "
$ws.Range("D25").Value = "'2.304"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +0.22%  "
$ws.Range("D26").Value = "2.128.85"
$ws.Range("E26").Value = "  -8.00%  "
$ws.Range("D27").Value = "'153.05"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -2.22%  "
$ws.Range("D28").Value = "'20.01"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -3.03%  "
$ws.Range("D29").Value = "'5.743"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -6.90%  "
$ws.Range("D30").Value = "'2.124"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -3.91%  "
$ws.Range("D31").Value = "'123.45"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -2.13%  "
$ws.Range("D32").Value = "'0.09640"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -1.66%  "
$ws.Range("D33").Value = "'0.9721"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -5.52%  "
$ws.Range("D34").Value = "'1.452"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -3.57%  "
$ws.Range("B35").Value = "HuobiToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D35").Value = "'3.636"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -2.93%  "
$ws.Range("B36").Value = "Filecoin"
$ws.Range("C36").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D36").Value = "'5.541"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -4.26%  "
$ws.Range("D37").Value = "'1.282"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -1.52%  "
$ws.Range("D38").Value = "'8.790"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -1.77%  "
$ws.Range("D39").Value = "'0.02291"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -4.01%  "
$ws.Range("E40").Value = "  -2.65%  "
$ws.Range("D41").Value = "'0.6147"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -4.33%  "
$ws.Range("D42").Value = "'10.89"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -4.42%  "
$ws.Range("E43").Value = "  -1.38%  "
$ws.Range("D44").Value = "'0.1904"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -2.72%  "
$ws.Range("E45").Value = "  -4.09%  "
$ws.Range("D46").Value = "'0.5855"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -4.47%  "
$ws.Range("D47").Value = "'12.82"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -2.64%  "
$ws.Range("D48").Value = "'2.013"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -5.66%  "
$ws.Range("D49").Value = "'3.434"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -1.51%  "
$ws.Range("D50").Value = "'0.06915"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -0.08%  "
$ws.Range("D51").Value = "'110.83"
$ws.Range("D51").ClearFormats()
